$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.142.16"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").Value = "3.165.90"
$ws.Range("E3").Value = "  -4.64%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'587.58"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").Value = "'134.39"
$ws.Range("E6").Value = "  -6.53%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.161.37"
$ws.Range("E8").Value = "  -4.76%  "
$ws.Range("D9").Value = "'0.500"
$ws.Range("E9").Value = "  -4.67%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("E10").Value = "  -5.94%  "
$ws.Range("D11").Value = "'5.22"
$ws.Range("E11").Value = "  -5.31%  "
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  -6.28%  "
$ws.Range("D13").Value = "'0.0000233"
$ws.Range("E13").Value = "  -7.06%  "
$ws.Range("D14").Value = "'33.03"
$ws.Range("E14").Value = "  -5.35%  "
$ws.Range("D15").Value = "3.700.55"
$ws.Range("E15").Value = "  -4.25%  "
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").Value = "3.177.02"
$ws.Range("E17").Value = "  -4.18%  "
$ws.Range("D18").Value = "62.165.34"
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("D19").Value = "'6.53"
$ws.Range("E19").Value = "  -5.79%  "
$ws.Range("D20").Value = "'455.49"
$ws.Range("E20").Value = "  -5.97%  "
$ws.Range("D21").Value = "'13.82"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").Value = "'0.699"
$ws.Range("E22").Value = "  -5.84%  "
$ws.Range("D23").Value = "'7.57"
$ws.Range("E23").Value = "  -5.87%  "
$ws.Range("D24").Value = "'13.21"
$ws.Range("E24").Value = "  -3.10%  "
$ws.Range("D25").Value = "'82.08"
$ws.Range("E25").Value = "  -3.25%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("B27").Value = "FirstDigitalUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.67"
$ws.Range("E28").Value = "  -4.08%  "
$ws.Range("D29").Value = "'6.87"
$ws.Range("E29").Value = "  -5.85%  "
$ws.Range("D30").Value = "'7.77"
$ws.Range("E30").Value = "  -5.99%  "
$ws.Range("D31").Value = "'2.02"
$ws.Range("E31").Value = "  -6.78%  "
$ws.Range("D32").Value = "'27.08"
$ws.Range("E32").Value = "  -8.46%  "
$ws.Range("E33").Value = "  -4.83%  "
$ws.Range("D34").Value = "'2.38"
$ws.Range("E34").Value = "  -7.51%  "
$ws.Range("D35").Value = "'1.03"
$ws.Range("E35").Value = "  -6.90%  "
$ws.Range("D36").Value = "'5.76"
$ws.Range("E36").Value = "  -4.50%  "
$ws.Range("D37").Value = "'51.03"
$ws.Range("E37").Value = "  -4.37%  "
$ws.Range("D38").Value = "0.0₃0682"
$ws.Range("E38").Value = "  -10.37%  "
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("D40").Value = "2.942.21"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("D41").Value = "'407.96"
$ws.Range("E41").Value = "  -5.97%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.63"
$ws.Range("E42").Value = "  -5.89%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.112"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").Value = "'7.96"
$ws.Range("E44").Value = "  -5.94%  "
$ws.Range("D45").Value = "'0.248"
$ws.Range("E45").Value = "  -7.74%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.12"
$ws.Range("E47").Value = "  -4.72%  "
$ws.Range("D48").Value = "'35.55"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("D49").Value = "'25.32"
$ws.Range("E49").Value = "  -5.06%  "
$ws.Range("D50").Value = "'123.38"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  -4.38%  "
